$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 47

# Text columns: Date, Time, Weekday, Week.
# Force text number format before assignment so values that look like
# dates/numbers ("2024-01-11", "01") are stored as text, then clear the
# formatting again so no extra cell style is left behind.
$textCols = 1,2,3,4
foreach ($col in $textCols) {
    $ws.Cells.Item($row, $col).NumberFormat = "@"
}

$ws.Cells.Item($row, 1).Value = "2024-01-11"
$ws.Cells.Item($row, 2).Value = "14:49:01"
$ws.Cells.Item($row, 3).Value = "Thursday"
$ws.Cells.Item($row, 4).Value = "01"

foreach ($col in $textCols) {
    $ws.Cells.Item($row, $col).ClearFormats()
}

# Numeric columns: Beijing, Guangzhou, Suzhou, Hangzhou, Nanjing, Xi_an,
# Chengdu, Chongqing, Tianjin, Hefei, Fuzhou, Xiamen, Changsha, Shanghai,
# Shenzhen, Wuhan.
$ws.Cells.Item($row, 5).Value = 139484
$ws.Cells.Item($row, 6).Value = 142824
$ws.Cells.Item($row, 7).Value = 171590
$ws.Cells.Item($row, 8).Value = 148141
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 119253
$ws.Cells.Item($row, 11).Value = 224830
$ws.Cells.Item($row, 12).Value = 251896
$ws.Cells.Item($row, 13).Value = 185325
$ws.Cells.Item($row, 14).Value = 110442
$ws.Cells.Item($row, 15).Value = 40765
$ws.Cells.Item($row, 16).Value = 30883
$ws.Cells.Item($row, 17).Value = 72844
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42063
$ws.Cells.Item($row, 20).Value = -1
